# Add a new Product Backlog story (row 9) about an overall resource calendar.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Product Backlog")

$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Ressourcen Übersicht"
$ws.Cells.Item(9, 3).Value = "Auf einem Kalender sollen alle MA Einsätze dargestellt werden"
$ws.Cells.Item(9, 4).Value = "low"
$ws.Cells.Item(9, 5).Value = 25
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = "waiting"

# Match the style of the other "Status" cells (vertical-top alignment).
$ws.Range("H9").VerticalAlignment = -4160

# Restore the updated selection state for each sheet, as recorded after editing.
$ws.Range("C13").Select()

$sprintWs = $wb.Worksheets.Item("Sprint Backlog")
$sprintWs.Activate()
$sprintWs.Range("C19").Select()
$excel.ActiveWindow.ScrollRow = 7
